# Update the "SecurityGroupRule_List (2)" example rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 becomes the single remaining example row, with new values.
$ws.Range("B3").Value = "Outbound"
$ws.Range("C3").Value = "192.168.13.0/24"
$ws.Range("D3").Value = "직접입력(Direct Input)"
$ws.Range("E3").Value = "2866, 3306"

# Rows 4-7 keep their formatting but lose their sample data.
$ws.Range("B4:E7").ClearContents()

# Rows 8-13 are removed entirely (sheet shrinks from H13 to H7).
$ws.Rows("8:13").Delete()

# Move the active selection to C10, matching the saved view state.
$ws.Range("C10").Select()
